# Update Betfair back/lay odds for 2025-11-14 fixtures.
# Applies the updated odds values to the corresponding cells on the
# active worksheet, row by row (rows 2-10), matching the latest
# scrape for this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - SCR Altach vs SSV Ulm
$ws.Range("J2").Value = 3.7
$ws.Range("M2").Value = 1.02

# Row 3 - FC Gutersloh vs TSV Havelse
$ws.Range("F3").Value = 2
$ws.Range("J3").Value = 3.75
$ws.Range("M3").Value = 1.02
$ws.Range("V3").Value = 1.44
$ws.Range("W3").Value = 1.4

# Row 4 - FK Loznica vs Fk Smederevo
$ws.Range("F4").Value = 2.22
$ws.Range("G4").Value = 2.9
$ws.Range("H4").Value = 2.8
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 2.8
$ws.Range("K4").Value = 3.9
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 2.44
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 1.53
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.19
$ws.Range("S4").Value = 3.45
$ws.Range("T4").Value = 1.9
$ws.Range("U4").Value = 1.76
$ws.Range("V4").Value = 1.35
$ws.Range("W4").Value = 1.52

# Row 5 - HIK Hellerup vs Vendsyssel FF
$ws.Range("F5").Value = 3.25
$ws.Range("G5").Value = 3.9
$ws.Range("H5").Value = 2.14
$ws.Range("I5").Value = 2.44
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 3.9
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 3.45
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 1.86
$ws.Range("Q5").Value = 1.94
$ws.Range("R5").Value = 1.32
$ws.Range("S5").Value = 3.4
$ws.Range("T5").Value = 1.74
$ws.Range("U5").Value = 2.06
$ws.Range("V5").Value = 1.7
$ws.Range("W5").Value = 1.34
$ws.Range("AB5").Value = 16
$ws.Range("AD5").Value = 12
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 32
$ws.Range("AI5").Value = 1000
$ws.Range("AL5").Value = 65
$ws.Range("AO5").Value = 23

# Row 6 - Valladolid vs Las Palmas
$ws.Range("F6").Value = 2.4
$ws.Range("K6").Value = 3.2
$ws.Range("N6").Value = 2.74
$ws.Range("P6").Value = 1.57
$ws.Range("T6").Value = 2.14
$ws.Range("U6").Value = 1.81
$ws.Range("V6").Value = 1.36
$ws.Range("W6").Value = 1.68
$ws.Range("AB6").Value = 8.4

# Row 7 - Flint Town United vs The New Saints
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 1.26
$ws.Range("I7").Value = 1.34
$ws.Range("J7").Value = 6.4
$ws.Range("K7").Value = 7.8
$ws.Range("N7").Value = 7.8
$ws.Range("R7").Value = 2
$ws.Range("T7").Value = 1.72
$ws.Range("U7").Value = 2.1
$ws.Range("V7").Value = 3.9
$ws.Range("W7").Value = 1.08
$ws.Range("Z7").Value = 12.5
$ws.Range("AC7").Value = 20
$ws.Range("AD7").Value = 12.5
$ws.Range("AE7").Value = 14
$ws.Range("AF7").Value = 140
$ws.Range("AG7").Value = 44
$ws.Range("AH7").Value = 32
$ws.Range("AI7").Value = 34
$ws.Range("AJ7").Value = 400
$ws.Range("AK7").Value = 160
$ws.Range("AO7").Value = 3.45

# Row 8 - Cardiff Metropolitan vs Briton Ferry Llansawel
$ws.Range("F8").Value = 1.63
$ws.Range("G8").Value = 1.69
$ws.Range("H8").Value = 5.5
$ws.Range("I8").Value = 6.4
$ws.Range("J8").Value = 4.2
$ws.Range("K8").Value = 4.9
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 4.4
$ws.Range("O8").Value = 1.23
$ws.Range("P8").Value = 2.2
$ws.Range("Q8").Value = 1.67
$ws.Range("R8").Value = 1.48
$ws.Range("S8").Value = 2.7
$ws.Range("T8").Value = 1.75
$ws.Range("U8").Value = 2.08
$ws.Range("V8").Value = 1.18
$ws.Range("W8").Value = 2.44
$ws.Range("X8").Value = 21
$ws.Range("Y8").Value = 980
$ws.Range("Z8").Value = 55
$ws.Range("AA8").Value = 160
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 24
$ws.Range("AE8").Value = 80
$ws.Range("AF8").Value = 11.5
$ws.Range("AG8").Value = 10.5
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 75
$ws.Range("AJ8").Value = 16.5
$ws.Range("AK8").Value = 16.5
$ws.Range("AM8").Value = 140
$ws.Range("AN8").Value = 8.199999999999999
$ws.Range("AO8").Value = 80

# Row 9 - Lanus vs Atl Tucuman
$ws.Range("F9").Value = 1.88
$ws.Range("G9").Value = 2.04
$ws.Range("H9").Value = 4.5
$ws.Range("I9").Value = 5.2
$ws.Range("J9").Value = 3.35
$ws.Range("K9").Value = 3.6

# Row 10 - Paysandu vs Amazonas FC
$ws.Range("F10").Value = 2.86
$ws.Range("G10").Value = 3.15
$ws.Range("H10").Value = 2.62
$ws.Range("I10").Value = 2.84
$ws.Range("J10").Value = 3.15
$ws.Range("K10").Value = 3.55
$ws.Range("L10").Value = 1.47
$ws.Range("V10").Value = 1.55
$ws.Range("W10").Value = 1.46
$ws.Range("AA10").Value = 50
$ws.Range("AB10").Value = 10.5
$ws.Range("AG10").Value = 14
$ws.Range("AI10").Value = 55
$ws.Range("AJ10").Value = 55
